$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45035
$ws.Range("I2").Value = 'Primera'
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 19000
$ws.Range("M2").Value = 19500
$ws.Range("P2").Value = 1083
$ws.Range("D3").Value = 44533
$ws.Range("H3").Value = 'Cultivar XV región'
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 6000
$ws.Range("L3").Value = 7000
$ws.Range("M3").Value = 6500
$ws.Range("N3").Value = '$/caja 10 kilos'
$ws.Range("O3").Value = 'Región de Arica y Parinacota'
$ws.Range("P3").Value = 650
$ws.Range("Q3").Value = 10
$ws.Range("D4").Value = 44533
$ws.Range("I4").Value = 'Segunda'
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 4000
$ws.Range("L4").Value = 5000
$ws.Range("M4").Value = 4500
$ws.Range("P4").Value = 450
$ws.Range("D5").Value = 44783
$ws.Range("J5").Value = 150
$ws.Range("D6").Value = 44363
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 140
$ws.Range("K6").Value = 14000
$ws.Range("L6").Value = 15000
$ws.Range("M6").Value = 14500
$ws.Range("P6").Value = 806
$ws.Range("D7").Value = 45043
$ws.Range("H7").Value = 'Cultivar IV Región'
$ws.Range("J7").Value = 170
$ws.Range("K7").Value = 18000
$ws.Range("L7").Value = 20000
$ws.Range("M7").Value = 19059
$ws.Range("N7").Value = '$/bandeja 18 kilos'
$ws.Range("O7").Value = 'Provincia de Limarí'
$ws.Range("P7").Value = 1059
$ws.Range("Q7").Value = 18
$ws.Range("D8").Value = 44526
$ws.Range("J8").Value = 100
$ws.Range("L8").Value = 5500
$ws.Range("M8").Value = 5250
$ws.Range("P8").Value = 525
$ws.Range("D9").Value = 44526
$ws.Range("H9").Value = 'Cultivar XV región'
$ws.Range("I9").Value = 'Segunda'
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 4000
$ws.Range("L9").Value = 4500
$ws.Range("M9").Value = 4250
$ws.Range("N9").Value = '$/caja 10 kilos'
$ws.Range("O9").Value = 'Región de Arica y Parinacota'
$ws.Range("P9").Value = 425
$ws.Range("Q9").Value = 10
$ws.Range("D10").Value = 44526
$ws.Range("H10").Value = 'Cultivar XV región'
$ws.Range("I10").Value = 'Tercera'
$ws.Range("J10").Value = 120
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3500
$ws.Range("M10").Value = 3250
$ws.Range("N10").Value = '$/caja 10 kilos'
$ws.Range("O10").Value = 'Región de Arica y Parinacota'
$ws.Range("P10").Value = 325
$ws.Range("Q10").Value = 10
$ws.Range("D11").Value = 44776
$ws.Range("H11").Value = 'Cultivar IV Región'
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 17000
$ws.Range("L11").Value = 18000
$ws.Range("M11").Value = 17500
$ws.Range("N11").Value = '$/bandeja 18 kilos'
$ws.Range("O11").Value = 'Provincia de Limarí'
$ws.Range("P11").Value = 972
$ws.Range("Q11").Value = 18
$ws.Range("D12").Value = 44377
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 100
$ws.Range("M12").Value = 17600
$ws.Range("P12").Value = 978
$ws.Range("D14").Value = 45042
$ws.Range("I14").Value = 'Segunda'
$ws.Range("J14").Value = 220
$ws.Range("M14").Value = 17545
$ws.Range("P14").Value = 975
$ws.Range("D15").Value = 44454
$ws.Range("J15").Value = 160
$ws.Range("K15").Value = 19000
$ws.Range("L15").Value = 20000
$ws.Range("M15").Value = 19500
$ws.Range("P15").Value = 1083
$ws.Range("D16").Value = 44398
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 17000
$ws.Range("L16").Value = 18000
$ws.Range("M16").Value = 17500
$ws.Range("P16").Value = 972
$ws.Range("D17").Value = 44398
$ws.Range("H17").Value = 'Cultivar IV Región'
$ws.Range("I17").Value = 'Segunda'
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 16000
$ws.Range("M17").Value = 15500
$ws.Range("N17").Value = '$/bandeja 18 kilos'
$ws.Range("O17").Value = 'Provincia de Limarí'
$ws.Range("P17").Value = 861
$ws.Range("Q17").Value = 18
$ws.Range("D18").Value = 44762
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 15000
$ws.Range("L18").Value = 16000
$ws.Range("M18").Value = 15500
$ws.Range("P18").Value = 861
$ws.Range("D19").Value = 44405
$ws.Range("I19").Value = 'Segunda'
$ws.Range("J19").Value = 140
$ws.Range("K19").Value = 17000
$ws.Range("L19").Value = 18000
$ws.Range("M19").Value = 17500
$ws.Range("P19").Value = 972
$ws.Range("D20").Value = 44755
$ws.Range("K20").Value = 17000
$ws.Range("L20").Value = 18000
$ws.Range("M20").Value = 17500
$ws.Range("P20").Value = 972
$ws.Range("D21").Value = 44748
$ws.Range("I21").Value = 'Primera'
$ws.Range("J21").Value = 250
$ws.Range("I22").Value = 'Segunda'
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 17000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 17500
$ws.Range("P22").Value = 972
$ws.Range("D23").Value = 44433
$ws.Range("I23").Value = 'Tercera'
$ws.Range("J23").Value = 120
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("P23").Value = 806
$ws.Range("D24").Value = 44435
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 100
$ws.Range("D25").Value = 44435
$ws.Range("I25").Value = 'Tercera'
$ws.Range("J25").Value = 120
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("P25").Value = 806
$ws.Range("D26").Value = 44412
$ws.Range("H26").Value = 'Cultivar IV Región'
$ws.Range("J26").Value = 150
$ws.Range("K26").Value = 17000
$ws.Range("L26").Value = 18000
$ws.Range("M26").Value = 17500
$ws.Range("N26").Value = '$/bandeja 18 kilos'
$ws.Range("O26").Value = 'Provincia de Limarí'
$ws.Range("P26").Value = 972
$ws.Range("Q26").Value = 18
$ws.Range("D27").Value = 44769
$ws.Range("H27").Value = 'Cultivar IV Región'
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 140
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17500
$ws.Range("N27").Value = '$/bandeja 18 kilos'
$ws.Range("O27").Value = 'Provincia de Limarí'
$ws.Range("P27").Value = 972
$ws.Range("Q27").Value = 18
$ws.Range("D28").Value = 44211
$ws.Range("I28").Value = 'Segunda'
$ws.Range("J28").Value = 140
$ws.Range("K28").Value = 4500
$ws.Range("L28").Value = 5000
$ws.Range("M28").Value = 4750
$ws.Range("P28").Value = 475
$ws.Range("D29").Value = 44757
$ws.Range("H29").Value = 'Cultivar XV región'
$ws.Range("I29").Value = 'Primera'
$ws.Range("J29").Value = 150
$ws.Range("K29").Value = 6000
$ws.Range("L29").Value = 6500
$ws.Range("M29").Value = 6250
$ws.Range("N29").Value = '$/caja 10 kilos'
$ws.Range("O29").Value = 'Región de Arica y Parinacota'
$ws.Range("P29").Value = 625
$ws.Range("Q29").Value = 10
$ws.Range("D30").Value = 44742
$ws.Range("J30").Value = 250
$ws.Range("D31").Value = 44554
$ws.Range("H31").Value = 'Cultivar XV región'
$ws.Range("J31").Value = 200
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 6000
$ws.Range("M31").Value = 5500
$ws.Range("N31").Value = '$/caja 10 kilos'
$ws.Range("O31").Value = 'Región de Arica y Parinacota'
$ws.Range("P31").Value = 550
$ws.Range("Q31").Value = 10
$ws.Range("D32").Value = 44771
$ws.Range("J32").Value = 140
$ws.Range("K32").Value = 8000
$ws.Range("L32").Value = 9000
$ws.Range("M32").Value = 8500
$ws.Range("P32").Value = 850
$ws.Range("D33").Value = 44221
$ws.Range("I33").Value = 'Primera'
$ws.Range("J33").Value = 140
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = 5500
$ws.Range("P33").Value = 550
$ws.Range("D34").Value = 44391
$ws.Range("I34").Value = 'Segunda'
$ws.Range("J34").Value = 100
$ws.Range("K34").Value = 15000
$ws.Range("L34").Value = 16000
$ws.Range("M34").Value = 15500
$ws.Range("P34").Value = 861
